# Adds a new "2022-Q4" sheet (with fund holdings data) right after the
# "总计" summary sheet, and updates the summary sheet with the new
# quarter's totals. All subsequent quarter sheets shift right
# automatically since they are addressed by name, not position.

function Set-TextValue {
    param($range, $value)
    # Force the cell to be stored as literal text (keeps leading zeros /
    # numeric-looking strings like "16.39" from being coerced into a
    # number) - the leading apostrophe is Excel's classic "treat as text"
    # quote-prefix marker and is not part of the stored value.
    $range.Value = "'" + $value
    # Quote-prefixing leaves a style behind (quotePrefix flag); resetting
    # to the Normal style afterwards clears it again while the cell keeps
    # its string type, matching the plain (un-styled) text cells in the
    # source data.
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# Match the page margins used by the other quarter sheets (0.75in /
# 1in / 0.5in) instead of the blank-sheet defaults (0.7in / 0.75in / 0.3in).
$ps = $newSheet.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# Borrow the cell formatting (header style + index-column style) from the
# "2021-Q4" sheet, which already has 10 rows - enough to cover our 12
# (1 header + 11 data) rows after extending by two more.
$fmtSource = $wb.Worksheets.Item("2021-Q4")
$fmtSource.Range("A1:H10").Copy()
$newSheet.Range("A1:H10").PasteSpecial(-4122)
$fmtSource.Range("A10:H10").Copy()
$newSheet.Range("A11:H12").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Populate the header row.
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $headers.Count; $col++) {
    $newSheet.Cells.Item(1, $col + 2).Value = $headers[$col]
}

# ---------------------------------------------------------------------
# 3. Populate the fund rows.
#    columns: code, name, scale, totalPosition, positionPct, marketValue, rank
# ---------------------------------------------------------------------
$rows = @(
    @("001643", "汇丰晋信智造先锋股票A", "16.39", "93.66", "3.62", "0.5933", 9),
    @("001644", "汇丰晋信智造先锋股票C", "8.77", "93.66", "3.62", "0.3175", 9),
    @("012850", "中融低碳经济3个月持有期混合A", "5.54", "88.38", "2.94", "0.1629", 9),
    @("006234", "万家汽车新趋势混合C", "4.06", "90.27", "3.41", "0.1384", 8),
    @("014575", "鑫元清洁能源混合C", "1.28", "94.25", "6.51", "0.0833", 9),
    @("001742", "广发百发大数据策略精选灵活配置混合E", "2.01", "51.25", "3.26", "0.0655", 8),
    @("006233", "万家汽车新趋势混合A", "1.65", "90.27", "3.41", "0.0563", 8),
    @("014574", "鑫元清洁能源混合A", "0.66", "94.25", "6.51", "0.0430", 9),
    @("012851", "中融低碳经济3个月持有期混合C", "1.19", "88.38", "2.94", "0.0350", 9),
    @("010487", "中银顺盈回报一年持有期混合", "0.75", "21.31", "1.16", "0.0087", 2),
    @("001741", "广发百发大数据策略精选灵活配置混合A", "0.03", "51.25", "3.26", "0.0010", 8)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    Set-TextValue $newSheet.Cells.Item($r, 2) $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    Set-TextValue $newSheet.Cells.Item($r, 4) $row[2]
    Set-TextValue $newSheet.Cells.Item($r, 5) $row[3]
    Set-TextValue $newSheet.Cells.Item($r, 6) $row[4]
    Set-TextValue $newSheet.Cells.Item($r, 7) $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 4. Update the "总计" summary sheet: shift the existing quarters down
#    one row and insert the new 2022-Q4 totals at the top.
# ---------------------------------------------------------------------

# Copy the format of the last existing data row (row 6, "2021-Q3") down
# to the new row 7 so its index-column cell picks up the same style.
$totalSheet.Range("A6:D6").Copy()
$totalSheet.Range("A7:D7").PasteSpecial(-4122)

# Shift quarter rows down by one (bottom-up, so we don't clobber data we
# still need to read).
for ($r = 6; $r -ge 2; $r--) {
    $totalSheet.Cells.Item($r + 1, 2).Value = $totalSheet.Cells.Item($r, 2).Value()
    $totalSheet.Cells.Item($r + 1, 3).Value = $totalSheet.Cells.Item($r, 3).Value()
    $totalSheet.Cells.Item($r + 1, 4).Value = $totalSheet.Cells.Item($r, 4).Value()
}

# Index column stays the plain sequential 0..5.
$totalSheet.Cells.Item(7, 1).Value = 5

# New top row: the 2022-Q4 totals.
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 11
$totalSheet.Cells.Item(2, 4).Value = 1.5

# ---------------------------------------------------------------------
# 5. Restore the originally-selected tab ("2021-Q3", the last sheet) so
#    the workbook/sheet view state is unaffected by the new sheet insert.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
